{"js": "const replacements = [\n  [\"97\u00d757=5529\", \"48\u00d783=3984\"],\n  [\"27\u00d727=729\", \"93\u00d778=7254\"],\n  [\"75\u00d790=6750\", \"51\u00d776=3876\"],\n  [\"48\u00d755=2640\", \"64\u00d725=1600\"],\n  [\"13\u00d789=1157\", \"33\u00d745=1485\"],\n  [\"95\u00d762=5890\", \"81\u00d732=2592\"],\n  [\"71\u00d746=3266\", \"97\u00d711=1067\"],\n  [\"68\u00d731=2108\", \"97\u00d778=7566\"],\n  [\"94\u00d799=9306\", \"89\u00d785=7565\"],\n  [\"15\u00d787=1305\", \"13\u00d790=1170\"],\n  [\"96\u00d766=6336\", \"29\u00d741=1189\"],\n  [\"27\u00d735=945\", \"66\u00d747=3102\"],\n  [\"11\u00d786=946\", \"14\u00d788=1232\"],\n  [\"48\u00d761=2928\", \"21\u00d793=1953\"],\n  [\"69\u00d793=6417\", \"37\u00d711=407\"],\n  [\"42\u00d764=2688\", \"97\u00d754=5238\"],\n  [\"81\u00d782=6642\", \"86\u00d760=5160\"],\n  [\"39\u00d744=1716\", \"77\u00d779=6083\"],\n  [\"98\u00d763=6174\", \"71\u00d790=6390\"],\n  [\"19\u00d717=323\", \"37\u00d793=3441\"],\n  [\"75\u00d717=1275\", \"82\u00d774=6068\"],\n  [\"88\u00d722=1936\", \"42\u00d754=2268\"],\n  [\"23\u00d719=437\", \"76\u00d777=5852\"],\n  [\"33\u00d773=2409\", \"99\u00d752=5148\"],\n  [\"81\u00d765=5265\", \"74\u00d762=4588\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (const table of tables.items) {\n  const range = table.getRange();\n  for (const [oldText, newText] of replacements) {\n    const results = range.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n    for (const r of results.items) {\n      r.insertText(newText, Word.InsertLocation.replace);\n    }\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"97\u00d757=5529\", \"48\u00d783=3984\"),\n    @(\"27\u00d727=729\", \"93\u00d778=7254\"),\n    @(\"75\u00d790=6750\", \"51\u00d776=3876\"),\n    @(\"48\u00d755=2640\", \"64\u00d725=1600\"),\n    @(\"13\u00d789=1157\", \"33\u00d745=1485\"),\n    @(\"95\u00d762=5890\", \"81\u00d732=2592\"),\n    @(\"71\u00d746=3266\", \"97\u00d711=1067\"),\n    @(\"68\u00d731=2108\", \"97\u00d778=7566\"),\n    @(\"94\u00d799=9306\", \"89\u00d785=7565\"),\n    @(\"15\u00d787=1305\", \"13\u00d790=1170\"),\n    @(\"96\u00d766=6336\", \"29\u00d741=1189\"),\n    @(\"27\u00d735=945\", \"66\u00d747=3102\"),\n    @(\"11\u00d786=946\", \"14\u00d788=1232\"),\n    @(\"48\u00d761=2928\", \"21\u00d793=1953\"),\n    @(\"69\u00d793=6417\", \"37\u00d711=407\"),\n    @(\"42\u00d764=2688\", \"97\u00d754=5238\"),\n    @(\"81\u00d782=6642\", \"86\u00d760=5160\"),\n    @(\"39\u00d744=1716\", \"77\u00d779=6083\"),\n    @(\"98\u00d763=6174\", \"71\u00d790=6390\"),\n    @(\"19\u00d717=323\", \"37\u00d793=3441\"),\n    @(\"75\u00d717=1275\", \"82\u00d774=6068\"),\n    @(\"88\u00d722=1936\", \"42\u00d754=2268\"),\n    @(\"23\u00d719=437\", \"76\u00d777=5852\"),\n    @(\"33\u00d773=2409\", \"99\u00d752=5148\"),\n    @(\"81\u00d765=5265\", \"74\u00d762=4588\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute(\n        $oldText,\n        $true,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $newText,\n        2\n    ) | Out-Null\n}\n\n$d.Save()\n"}
